$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 524.8333  # H2 was 449.77777
$ws.Cells.Item(2, 9).Value = 530  # I2 was 458.33334
$ws.Cells.Item(2, 10).Value = 499  # J2 was 432.66666
$ws.Cells.Item(2, 11).Value = 530  # K2 was 458.33334
$ws.Cells.Item(2, 12).Value = 499  # L2 was 432.66666
$ws.Cells.Item(2, 13).Value = -417  # M2 was -345.33334
$ws.Cells.Item(2, 14).Value = -725  # N2 was -658.66666
$ws.Cells.Item(12, 8).Value = 427  # H12 was 389.55554
$ws.Cells.Item(12, 10).Value = 200  # J12 was 163.33333
$ws.Cells.Item(12, 12).Value = 200  # L12 was 163.33333
$ws.Cells.Item(12, 14).Value = -540  # N12 was -503.33333
$ws.Cells.Item(40, 8).Value = 836706  # H40 was 772444
$ws.Cells.Item(40, 10).Value = 2443  # J40 was 2214.4
$ws.Cells.Item(40, 12).Value = 2443  # L40 was 2214.4
$ws.Cells.Item(40, 14).Value = -2793  # N40 was -2564.4
$ws.Cells.Item(62, 8).Value = 8806.538  # H62 was 8713.214
$ws.Cells.Item(62, 10).Value = 8915.583000000001  # J62 was 8806.691999999999
$ws.Cells.Item(62, 12).Value = 8915.583000000001  # L62 was 8806.691999999999
$ws.Cells.Item(62, 14).Value = -10163.583  # N62 was -10054.692
$ws.Cells.Item(65, 8).Value = 8806.538  # H65 was 8713.214
$ws.Cells.Item(65, 10).Value = 8915.583000000001  # J65 was 8806.691999999999
$ws.Cells.Item(65, 12).Value = 44577.915  # L65 was 44033.45999999999
$ws.Cells.Item(65, 14).Value = -50817.915  # N65 was -50273.45999999999
$ws.Cells.Item(92, 8).Value = 1727.0769  # H92 was 1646.0834
$ws.Cells.Item(92, 9).Value = 1364.8  # I92 was 1394.8
$ws.Cells.Item(92, 10).Value = 2934.6667  # J92 was 2902.5
$ws.Cells.Item(92, 11).Value = 1364.8  # K92 was 1394.8
$ws.Cells.Item(92, 12).Value = 2934.6667  # L92 was 2902.5
$ws.Cells.Item(92, 13).Value = -116.8  # M92 was -146.8
$ws.Cells.Item(92, 14).Value = -5430.6667  # N92 was -5398.5
$ws.Cells.Item(116, 8).Value = 6050  # H116 was 5899.8335
$ws.Cells.Item(116, 9).Value = 6366.6665  # I116 was 6075
$ws.Cells.Item(116, 10).Value = 5100  # J116 was 5549.5
$ws.Cells.Item(116, 11).Value = 6366.6665  # K116 was 6075
$ws.Cells.Item(116, 12).Value = 5100  # L116 was 5549.5
$ws.Cells.Item(116, 13).Value = -2924.6665  # M116 was -2633
$ws.Cells.Item(116, 14).Value = -11984  # N116 was -12433.5
$ws.Cells.Item(137, 8).Value = 3466.5557  # H137 was 3639.88
$ws.Cells.Item(137, 10).Value = 4471.5625  # J137 was 4924.643
$ws.Cells.Item(137, 12).Value = 13414.6875  # L137 was 14773.929
$ws.Cells.Item(137, 14).Value = -18514.6875  # N137 was -19873.929
$ws.Cells.Item(138, 8).Value = 6814.939  # H138 was 6768.64
$ws.Cells.Item(138, 10).Value = 7888.6924  # J138 was 7803.975
$ws.Cells.Item(138, 12).Value = 23666.0772  # L138 was 23411.925
$ws.Cells.Item(138, 14).Value = -33946.0772  # N138 was -33691.925

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2244.6667  # H32 was 2267.3538
$ws.Cells.Item(32, 9).Value = 1330.9166  # I32 was 1368.9178
$ws.Cells.Item(32, 11).Value = 1330.9166  # K32 was 1368.9178
$ws.Cells.Item(32, 13).Value = -1043.9166  # M32 was -1081.9178
$ws.Cells.Item(38, 8).Value = 50000  # H38 was 0
$ws.Cells.Item(38, 9).Value = 50000  # I38 was 0
$ws.Cells.Item(38, 11).Value = 50000  # K38 was 0
$ws.Cells.Item(38, 13).Value = -49533  # M38 was None
$ws.Cells.Item(74, 8).Value = 20857988  # H74 was 21301732
$ws.Cells.Item(74, 10).Value = 2641.5715  # J74 was 2740.1667
$ws.Cells.Item(74, 12).Value = 2641.5715  # L74 was 2740.1667
$ws.Cells.Item(74, 14).Value = -4389.5715  # N74 was -4488.1667
$ws.Cells.Item(77, 8).Value = 20857988  # H77 was 21301732
$ws.Cells.Item(77, 10).Value = 2641.5715  # J77 was 2740.1667
$ws.Cells.Item(77, 12).Value = 13207.8575  # L77 was 13700.8335
$ws.Cells.Item(77, 14).Value = -21943.8575  # N77 was -22436.8335

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3164.7368  # H99 was 3530
$ws.Cells.Item(99, 9).Value = 2399.8333  # I99 was 2794.2222
$ws.Cells.Item(99, 11).Value = 2399.8333  # K99 was 2794.2222
$ws.Cells.Item(99, 13).Value = -901.8332999999998  # M99 was -1296.2222

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2295.9  # H16 was 2460.889
$ws.Cells.Item(16, 9).Value = 2183  # I16 was 2457.4
$ws.Cells.Item(16, 11).Value = 2183  # K16 was 2457.4
$ws.Cells.Item(16, 13).Value = -1896  # M16 was -2170.4
$ws.Cells.Item(31, 8).Value = 5003165.5  # H31 was 5558952.5
$ws.Cells.Item(31, 9).Value = 1936.9286  # I31 was 2138.4583
$ws.Cells.Item(31, 10).Value = 11368366  # J31 was 11909598
$ws.Cells.Item(31, 11).Value = 1936.9286  # K31 was 2138.4583
$ws.Cells.Item(31, 12).Value = 11368366  # L31 was 11909598
$ws.Cells.Item(31, 13).Value = -1641.9286  # M31 was -1843.4583
$ws.Cells.Item(31, 14).Value = -11368956  # N31 was -11910188
$ws.Cells.Item(34, 8).Value = 5003165.5  # H34 was 5558952.5
$ws.Cells.Item(34, 9).Value = 1936.9286  # I34 was 2138.4583
$ws.Cells.Item(34, 10).Value = 11368366  # J34 was 11909598
$ws.Cells.Item(34, 11).Value = 1936.9286  # K34 was 2138.4583
$ws.Cells.Item(34, 12).Value = 11368366  # L34 was 11909598
$ws.Cells.Item(34, 13).Value = -1734.9286  # M34 was -1936.4583
$ws.Cells.Item(34, 14).Value = -11368770  # N34 was -11910002
$ws.Cells.Item(41, 8).Value = 16420.334  # H41 was 21062.834
$ws.Cells.Item(41, 9).Value = 1381.125  # I41 was 2194.1428
$ws.Cells.Item(41, 10).Value = 46498.75  # J41 was 47479
$ws.Cells.Item(41, 11).Value = 1381.125  # K41 was 2194.1428
$ws.Cells.Item(41, 12).Value = 46498.75  # L41 was 47479
$ws.Cells.Item(41, 13).Value = -953.125  # M41 was -1766.1428
$ws.Cells.Item(41, 14).Value = -47354.75  # N41 was -48335
$ws.Cells.Item(94, 8).Value = 2415.5  # H94 was 2738.8
$ws.Cells.Item(94, 9).Value = 799.5  # I94 was 800
$ws.Cells.Item(94, 11).Value = 799.5  # K94 was 800
$ws.Cells.Item(94, 13).Value = -348.5  # M94 was -349
$ws.Cells.Item(99, 8).Value = 7731.7026  # H99 was 8556.807000000001
$ws.Cells.Item(99, 9).Value = 4268.846  # I99 was 4988.5557
$ws.Cells.Item(99, 10).Value = 9607.416999999999  # J99 was 10016.546
$ws.Cells.Item(99, 11).Value = 4268.846  # K99 was 4988.5557
$ws.Cells.Item(99, 12).Value = 9607.416999999999  # L99 was 10016.546
$ws.Cells.Item(99, 13).Value = -2770.846  # M99 was -3490.5557
$ws.Cells.Item(99, 14).Value = -12603.417  # N99 was -13012.546
$ws.Cells.Item(113, 8).Value = 2295.9  # H113 was 2460.889
$ws.Cells.Item(113, 9).Value = 2183  # I113 was 2457.4
$ws.Cells.Item(113, 11).Value = 2183  # K113 was 2457.4
$ws.Cells.Item(113, 13).Value = -13  # M113 was -287.4000000000001
$ws.Cells.Item(126, 8).Value = 7731.7026  # H126 was 8556.807000000001
$ws.Cells.Item(126, 9).Value = 4268.846  # I126 was 4988.5557
$ws.Cells.Item(126, 10).Value = 9607.416999999999  # J126 was 10016.546
$ws.Cells.Item(126, 11).Value = 12806.538  # K126 was 14965.6671
$ws.Cells.Item(126, 12).Value = 28822.251  # L126 was 30049.638
$ws.Cells.Item(126, 13).Value = -10336.538  # M126 was -12495.6671
$ws.Cells.Item(126, 14).Value = -33762.251  # N126 was -34989.638

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 1132  # H51 was 2797.2
$ws.Cells.Item(51, 9).Value = 1132  # I51 was 2496.5
$ws.Cells.Item(51, 10).Value = 0  # J51 was 4000
$ws.Cells.Item(51, 11).Value = 3396  # K51 was 7489.5
$ws.Cells.Item(51, 12).Value = 0  # L51 was 12000
$ws.Cells.Item(51, 13).Value = -2936  # M51 was -7029.5
$ws.Cells.Item(51, 14).ClearContents()  # N51 was -12920
$ws.Cells.Item(68, 8).Value = 2359.4062  # H68 was 2419.7097
$ws.Cells.Item(68, 9).Value = 1868.5333  # I68 was 1967
$ws.Cells.Item(68, 11).Value = 5605.5999  # K68 was 5901
$ws.Cells.Item(68, 13).Value = -4794.5999  # M68 was -5090
$ws.Cells.Item(71, 8).Value = 2359.4062  # H71 was 2419.7097
$ws.Cells.Item(71, 9).Value = 1868.5333  # I71 was 1967
$ws.Cells.Item(71, 11).Value = 16816.7997  # K71 was 17703
$ws.Cells.Item(71, 13).Value = -12760.7997  # M71 was -13647
$ws.Cells.Item(113, 8).Value = 855.7273  # H113 was 796.9167
$ws.Cells.Item(113, 9).Value = 181.25  # I113 was 177.77777
$ws.Cells.Item(113, 11).Value = 543.75  # K113 was 533.33331
$ws.Cells.Item(113, 13).Value = 1626.25  # M113 was 1636.66669
$ws.Cells.Item(122, 8).Value = 1372.6923  # H122 was 1351.2667
$ws.Cells.Item(122, 9).Value = 747.4286  # I122 was 832
$ws.Cells.Item(122, 10).Value = 2102.1667  # J122 was 1697.4445
$ws.Cells.Item(122, 11).Value = 6726.8574  # K122 was 7488
$ws.Cells.Item(122, 12).Value = 18919.5003  # L122 was 15277.0005
$ws.Cells.Item(122, 13).Value = -4276.8574  # M122 was -5038
$ws.Cells.Item(122, 14).Value = -23819.5003  # N122 was -20177.0005
$ws.Cells.Item(123, 8).Value = 3016.5  # H123 was 3033
$ws.Cells.Item(123, 9).Value = 3000  # I123 was 0
$ws.Cells.Item(123, 11).Value = 9000  # K123 was 0
$ws.Cells.Item(123, 13).Value = -6550  # M123 was None
$ws.Cells.Item(124, 8).Value = 3557.75  # H124 was 4126.2
$ws.Cells.Item(124, 9).Value = 2527.1428  # I124 was 2815
$ws.Cells.Item(124, 10).Value = 5000.6  # J124 was 5000.3335
$ws.Cells.Item(124, 11).Value = 7581.428400000001  # K124 was 8445
$ws.Cells.Item(124, 12).Value = 15001.8  # L124 was 15001.0005
$ws.Cells.Item(124, 13).Value = -2671.428400000001  # M124 was -3535
$ws.Cells.Item(124, 14).Value = -24821.8  # N124 was -24821.0005
$ws.Cells.Item(131, 8).Value = 41351.344  # H131 was 42679.07
$ws.Cells.Item(131, 10).Value = 9792.666999999999  # J131 was 10193.929
$ws.Cells.Item(131, 12).Value = 29378.001  # L131 was 30581.787
$ws.Cells.Item(131, 14).Value = -39458.001  # N131 was -40661.787
$ws.Cells.Item(132, 8).Value = 3564.7368  # H132 was 3580.7026
$ws.Cells.Item(132, 10).Value = 6138.4165  # J132 was 6426.091
$ws.Cells.Item(132, 12).Value = 55245.7485  # L132 was 57834.819
$ws.Cells.Item(132, 14).Value = -60305.7485  # N132 was -62894.819

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 6589.6  # H41 was 6991.3335
$ws.Cells.Item(126, 8).Value = 86924.75  # H126 was 86924.836
$ws.Cells.Item(126, 9).Value = 169016.17  # I126 was 169016.33
$ws.Cells.Item(126, 11).Value = 507048.51  # K126 was 507048.99
$ws.Cells.Item(126, 13).Value = -504578.51  # M126 was -504578.99

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3133.3845  # H40 was 3227.625
$ws.Cells.Item(40, 9).Value = 3122.182  # I40 was 3175.6191
$ws.Cells.Item(40, 10).Value = 3195  # J40 was 3591.6667
$ws.Cells.Item(40, 11).Value = 3122.182  # K40 was 3175.6191
$ws.Cells.Item(40, 12).Value = 3195  # L40 was 3591.6667
$ws.Cells.Item(40, 13).Value = -2986.182  # M40 was -3039.6191
$ws.Cells.Item(40, 14).Value = -3467  # N40 was -3863.6667
$ws.Cells.Item(61, 8).Value = 3032.3076  # H61 was 2882.9092
$ws.Cells.Item(61, 9).Value = 1894.6  # I61 was 1936.625
$ws.Cells.Item(61, 10).Value = 6824.6665  # J61 was 5406.3335
$ws.Cells.Item(61, 11).Value = 1894.6  # K61 was 1936.625
$ws.Cells.Item(61, 12).Value = 6824.6665  # L61 was 5406.3335
$ws.Cells.Item(61, 13).Value = -1692.6  # M61 was -1734.625
$ws.Cells.Item(61, 14).Value = -7228.6665  # N61 was -5810.3335
$ws.Cells.Item(113, 8).Value = 3032.3076  # H113 was 2882.9092
$ws.Cells.Item(113, 9).Value = 1894.6  # I113 was 1936.625
$ws.Cells.Item(113, 10).Value = 6824.6665  # J113 was 5406.3335
$ws.Cells.Item(113, 11).Value = 1894.6  # K113 was 1936.625
$ws.Cells.Item(113, 12).Value = 6824.6665  # L113 was 5406.3335
$ws.Cells.Item(113, 13).Value = 275.4000000000001  # M113 was 233.375
$ws.Cells.Item(113, 14).Value = -11164.6665  # N113 was -9746.333500000001
$ws.Cells.Item(123, 8).Value = 74998  # H123 was 30390
$ws.Cells.Item(123, 9).Value = 0  # I123 was 30390
$ws.Cells.Item(123, 10).Value = 74998  # J123 was 0
$ws.Cells.Item(123, 11).Value = 0  # K123 was 30390
$ws.Cells.Item(123, 12).Value = 74998  # L123 was 0
$ws.Cells.Item(123, 13).ClearContents()  # M123 was -25490
$ws.Cells.Item(123, 14).Value = -84798  # N123 was None

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 527271  # H100 was 500922.34
$ws.Cells.Item(100, 9).Value = 833995.3  # I100 was 769864.75
$ws.Cells.Item(100, 11).Value = 1667990.6  # K100 was 1539729.5
$ws.Cells.Item(100, 13).Value = -1667449.6  # M100 was -1539188.5
$ws.Cells.Item(107, 8).Value = 1665.5172  # H107 was 1762.931
$ws.Cells.Item(107, 9).Value = 1171.6471  # I107 was 1229.25
$ws.Cells.Item(107, 10).Value = 2365.1667  # J107 was 2419.7693
$ws.Cells.Item(107, 11).Value = 3514.9413  # K107 was 3687.75
$ws.Cells.Item(107, 12).Value = 7095.500100000001  # L107 was 7259.3079
$ws.Cells.Item(107, 13).Value = -1594.9413  # M107 was -1767.75
$ws.Cells.Item(107, 14).Value = -10935.5001  # N107 was -11099.3079
$ws.Cells.Item(113, 8).Value = 1245.591  # H113 was 1311.45
$ws.Cells.Item(113, 9).Value = 1024.4667  # I113 was 1056.6428
$ws.Cells.Item(113, 10).Value = 1719.4286  # J113 was 1906
$ws.Cells.Item(113, 11).Value = 3073.4001  # K113 was 3169.9284
$ws.Cells.Item(113, 12).Value = 5158.2858  # L113 was 5718
$ws.Cells.Item(113, 13).Value = -903.4000999999998  # M113 was -999.9284000000002
$ws.Cells.Item(113, 14).Value = -9498.2858  # N113 was -10058
$ws.Cells.Item(122, 8).Value = 7146443.5  # H122 was 7146526.5
$ws.Cells.Item(122, 9).Value = 2623.7368  # I122 was 2745.7896
$ws.Cells.Item(122, 11).Value = 7871.2104  # K122 was 8237.3688
$ws.Cells.Item(122, 13).Value = -5421.2104  # M122 was -5787.3688
$ws.Cells.Item(126, 8).Value = 12830798  # H126 was 16679826
$ws.Cells.Item(126, 9).Value = 15161953  # I126 was 18532042
$ws.Cells.Item(126, 10).Value = 9447.5  # J126 was 9895
$ws.Cells.Item(126, 11).Value = 45485859  # K126 was 55596126
$ws.Cells.Item(126, 12).Value = 28342.5  # L126 was 29685
$ws.Cells.Item(126, 13).Value = -45483389  # M126 was -55593656
$ws.Cells.Item(126, 14).Value = -33282.5  # N126 was -34625
$ws.Cells.Item(127, 8).Value = 199780  # H127 was 199828.75
$ws.Cells.Item(136, 8).Value = 3121.5715  # H136 was 3269.4614
$ws.Cells.Item(136, 9).Value = 3147.077  # I136 was 3309.4167
$ws.Cells.Item(136, 11).Value = 9441.231  # K136 was 9928.250100000001
$ws.Cells.Item(136, 13).Value = -6891.231  # M136 was -7378.250100000001
